# Add a new "2021" column (R) to the Sustainable Development Goal 3.9.2
# indicator table, mirroring the formatting already used by the
# neighbouring (N/O/P) year columns.
#
# Row 4  : year header (2021)
# Row 5  : first data row (bold style, like the "Kyrgyz Republic" total row)
# Rows 6-13: regular data rows
# Row 14 : last data row (bottom border)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the per-row formatting from column P (the last existing year
#     column before the totals gap) into column R, row by row, so the new
#     column gets exactly the same look (font/border/number format) as the
#     rest of the table. ---
$ws.Range("P4").Copy($ws.Range("R4"))
$ws.Range("P5").Copy($ws.Range("R5"))
$ws.Range("P6").Copy($ws.Range("R6"))
$ws.Range("P7").Copy($ws.Range("R7"))
$ws.Range("P8").Copy($ws.Range("R8"))
$ws.Range("P9").Copy($ws.Range("R9"))
$ws.Range("P10").Copy($ws.Range("R10"))
$ws.Range("P11").Copy($ws.Range("R11"))
$ws.Range("P12").Copy($ws.Range("R12"))
$ws.Range("P13").Copy($ws.Range("R13"))
$ws.Range("P14").Copy($ws.Range("R14"))

# --- Now write the real 2021 values/labels on top of the copied format. ---
$ws.Range("R4").Value = 2021

$ws.Range("R5").Value = 1
$ws.Range("R6").Value = 2.2
$ws.Range("R7").Value = 1.7
$ws.Range("R8").Value = "-"
$ws.Range("R9").Value = 0.3
$ws.Range("R10").Value = 1.1
$ws.Range("R11").Value = "-"
$ws.Range("R12").Value = 0.9
$ws.Range("R13").Value = 0.4
$ws.Range("R14").Value = 0.6

# --- Match the selection the author ended up with in the saved file. ---
$ws.Range("S17").Select()
